$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y5").NumberFormat = "@"

$ws.Range("A2").Value = 67753802
$ws.Range("B2").Value = 78569
$ws.Range("E2").Value = 6458
$ws.Range("F2").Value = "Lunglav"
$ws.Range("G2").Value = "Lobaria pulmonaria"
$ws.Range("H2").Value = "(L.) Hoffm."
$ws.Range("I2").Value = "2"
$ws.Range("P2").Value = "Hamrarna, Dlr"
$ws.Range("Q2").Value = 420563.8119189619
$ws.Range("R2").Value = 6725772.004727462
$ws.Range("Y2").Value = "2017-09-30"
$ws.Range("AA2").Value = "2017-09-30"
$ws.Range("AW2").Value = "per taube"
$ws.Range("AX2").Value = "per taube"
$ws.Range("A3").Value = 80491995
$ws.Range("Q3").Value = 420334.941858212
$ws.Range("R3").Value = 6725817.777769641
$ws.Range("A4").Value = 80491997
$ws.Range("Q4").Value = 420411.9141204953
$ws.Range("R4").Value = 6725811.167292034
$ws.Range("A5").Value = 80491992
$ws.Range("B5").Value = 77177
$ws.Range("E5").Value = 353
$ws.Range("F5").Value = "Dvärgbägarlav"
$ws.Range("G5").Value = "Cladonia parasitica"
$ws.Range("H5").Value = "(Hoffm.) Hoffm."
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("P5").Value = "Furubäcksmyrorna, Dlr"
$ws.Range("Q5").Value = 419982.7934210463
$ws.Range("R5").Value = 6725885.92551488
$ws.Range("Y5").Value = "2019-09-23"
$ws.Range("AA5").Value = "2019-09-23"
$ws.Range("AH5").Value = ""
$ws.Range("AJ5").Value = ""
$ws.Range("AK5").Value = ""
$ws.Range("AL5").Value = ""
$ws.Range("AO5").Value = ""
$ws.Range("AW5").Value = "Uno Skog"
$ws.Range("AX5").Value = "Uno Skog"
$ws.Range("A6").Value = 80492005
$ws.Range("B6").Value = 77506
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("P6").Value = "Hamrarna, Dlr"
$ws.Range("Q6").Value = 420638.7854734209
$ws.Range("R6").Value = 6725786.040621105
$ws.Range("A7").Value = 80491993
$ws.Range("B7").Value = 90653
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 4364
$ws.Range("F7").Value = "Dropptaggsvamp"
$ws.Range("G7").Value = "Hydnellum ferrugineum"
$ws.Range("H7").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P7").Value = "Furubäcksmyrorna, Dlr"
$ws.Range("Q7").Value = 420104.1243034533
$ws.Range("R7").Value = 6725864.098819775
$ws.Range("J2").Value = "dm²"
$ws.Range("K2").Value = ""
$ws.Range("AH2").Value = "Skogsmark"
$ws.Range("AJ2").Value = "asp"
$ws.Range("AK2").Value = "Populus tremula"
$ws.Range("AL2").Value = "ga asp m bohål, mkt fin naturskog på denna fastighet"
$ws.Range("AO2").Value = "Populus tremula # ga asp m bohål, mkt fin naturskog på denna fastighet"

$ws.Range("AA2").Style = "Normal"
$ws.Range("AA5").Style = "Normal"
$ws.Range("Y2").Style = "Normal"
$ws.Range("Y5").Style = "Normal"
